$d = $word.ActiveDocument
$s = $d.Styles.Add("TotallyNewStyleName", 1)
Write-Output $s.NameLocal
$s2 = $d.Styles.Add("egXML", 1)
Write-Output $s2.NameLocal
